$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 29 (shifts existing rows 29+ down by one) and populate it
# with the new image-compression data point.
$ws.Range("A29:E29").EntireRow.Insert()
$ws.Range("A29").Value = "data/BayOfBengal/Myanmar_mangrove_and_landuse_map.png"
$ws.Range("B29").Value = 13284671
$ws.Range("C29").Value = 3235628
$ws.Range("D29").Value = "jpg"
$ws.Range("E29").Formula = "=C29/B29"

# Keep the UI selection in sync with where Excel would have left it (one row
# further down than before, matching the new totals row).
$null = $ws.Range("B76").Select()
